# Scheduled runner update: refresh computed price/profit columns (H:N)
# for a set of leve rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW
# and WVR sheets, reflecting latest market board pricing data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 13.666667
$ws.Range("J6").Value = 20
$ws.Range("L6").Value = 60
$ws.Range("N6").Value = -284
$ws.Range("H33").Value = 377.11765
$ws.Range("I33").Value = 267.8889
$ws.Range("J33").Value = 500
$ws.Range("K33").Value = 267.8889
$ws.Range("L33").Value = 500
$ws.Range("M33").Value = -38.88889999999998
$ws.Range("N33").Value = -958
$ws.Range("H98").Value = 1901.5435
$ws.Range("I98").Value = 1896.9524
$ws.Range("K98").Value = 1896.9524
$ws.Range("M98").Value = -398.9523999999999
$ws.Range("H113").Value = 2827.6155
$ws.Range("I113").Value = 1976.25
$ws.Range("J113").Value = 3206
$ws.Range("K113").Value = 1976.25
$ws.Range("L113").Value = 3206
$ws.Range("M113").Value = 1277.75
$ws.Range("N113").Value = -9714
$ws.Range("H122").Value = 1901.5435
$ws.Range("I122").Value = 1896.9524
$ws.Range("K122").Value = 5690.857199999999
$ws.Range("M122").Value = -3240.857199999999
$ws.Range("H135").Value = 45456300
$ws.Range("I135").Value = 18520082
$ws.Range("K135").Value = 166680738
$ws.Range("M135").Value = -166678203
$ws.Range("H138").Value = 3783.9268
$ws.Range("J138").Value = 4734.6113
$ws.Range("L138").Value = 14203.8339
$ws.Range("N138").Value = -24483.8339

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8182.7183
$ws.Range("J32").Value = 19750
$ws.Range("L32").Value = 19750
$ws.Range("N32").Value = -20324
$ws.Range("H61").Value = 6382.298
$ws.Range("I61").Value = 3526.5938
$ws.Range("J61").Value = 12474.467
$ws.Range("K61").Value = 3526.5938
$ws.Range("L61").Value = 12474.467
$ws.Range("M61").Value = -3314.5938
$ws.Range("N61").Value = -12898.467
$ws.Range("H74").Value = 92663.25
$ws.Range("I74").Value = 102647.97
$ws.Range("K74").Value = 102647.97
$ws.Range("M74").Value = -101773.97
$ws.Range("H77").Value = 92663.25
$ws.Range("I77").Value = 102647.97
$ws.Range("K77").Value = 513239.85
$ws.Range("M77").Value = -508871.85
$ws.Range("H96").Value = 14936
$ws.Range("J96").Value = 14936
$ws.Range("L96").Value = 14936
$ws.Range("N96").Value = -20428
$ws.Range("H132").Value = 2698.6978
$ws.Range("I132").Value = 1576.5834
$ws.Range("K132").Value = 4729.7502
$ws.Range("M132").Value = -2199.7502
$ws.Range("H136").Value = 6382.298
$ws.Range("I136").Value = 3526.5938
$ws.Range("J136").Value = 12474.467
$ws.Range("K136").Value = 10579.7814
$ws.Range("L136").Value = 37423.401
$ws.Range("M136").Value = -8029.7814
$ws.Range("N136").Value = -42523.401

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 7097.8
$ws.Range("I7").Value = 530
$ws.Range("J7").Value = 16949.5
$ws.Range("K7").Value = 530
$ws.Range("L7").Value = 16949.5
$ws.Range("M7").Value = -417
$ws.Range("N7").Value = -17175.5
$ws.Range("H134").Value = 3320.1724
$ws.Range("J134").Value = 2842
$ws.Range("L134").Value = 8526
$ws.Range("N134").Value = -13596

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3768.3333
$ws.Range("I31").Value = 2075.3618
$ws.Range("J31").Value = 8741.4375
$ws.Range("K31").Value = 2075.3618
$ws.Range("L31").Value = 8741.4375
$ws.Range("M31").Value = -1780.3618
$ws.Range("N31").Value = -9331.4375
$ws.Range("H34").Value = 3768.3333
$ws.Range("I34").Value = 2075.3618
$ws.Range("J34").Value = 8741.4375
$ws.Range("K34").Value = 2075.3618
$ws.Range("L34").Value = 8741.4375
$ws.Range("M34").Value = -1873.3618
$ws.Range("N34").Value = -9145.4375
$ws.Range("H58").Value = 2167268.2
$ws.Range("I58").Value = 3789915.8
$ws.Range("J58").Value = 3738.2778
$ws.Range("K58").Value = 3789915.8
$ws.Range("L58").Value = 3738.2778
$ws.Range("M58").Value = -3789712.8
$ws.Range("N58").Value = -4144.2778
$ws.Range("H132").Value = 2881.5
$ws.Range("I132").Value = 2411.9375
$ws.Range("J132").Value = 3632.8
$ws.Range("K132").Value = 7235.8125
$ws.Range("L132").Value = 10898.4
$ws.Range("M132").Value = -4705.8125
$ws.Range("N132").Value = -15958.4
$ws.Range("H136").Value = 2167268.2
$ws.Range("I136").Value = 3789915.8
$ws.Range("J136").Value = 3738.2778
$ws.Range("K136").Value = 11369747.4
$ws.Range("L136").Value = 11214.8334
$ws.Range("M136").Value = -11367197.4
$ws.Range("N136").Value = -16314.8334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 81.666664
$ws.Range("I6").Value = 81.666664
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 244.999992
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -131.999992
$ws.Range("N6").Value = $null
$ws.Range("H113").Value = 806.3617
$ws.Range("I113").Value = 811.0833
$ws.Range("J113").Value = 790.9091
$ws.Range("K113").Value = 2433.2499
$ws.Range("L113").Value = 2372.7273
$ws.Range("M113").Value = -263.2498999999998
$ws.Range("N113").Value = -6712.7273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 10000
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 10000
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = $null
$ws.Range("N11").Value = -10278
$ws.Range("H19").Value = 10000
$ws.Range("J19").Value = 10000
$ws.Range("L19").Value = 10000
$ws.Range("N19").Value = -10576
$ws.Range("H102").Value = 4421.2256
$ws.Range("I102").Value = 4088.5908
$ws.Range("J102").Value = 5234.3335
$ws.Range("K102").Value = 4088.5908
$ws.Range("L102").Value = 5234.3335
$ws.Range("M102").Value = -2466.5908
$ws.Range("N102").Value = -8478.333500000001
$ws.Range("H132").Value = 25156.795
$ws.Range("I132").Value = 57033.39
$ws.Range("K132").Value = 171100.17
$ws.Range("M132").Value = -168570.17

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 24874.846
$ws.Range("I61").Value = 31947.3
$ws.Range("K61").Value = 31947.3
$ws.Range("M61").Value = -31745.3
$ws.Range("H113").Value = 24874.846
$ws.Range("I113").Value = 31947.3
$ws.Range("K113").Value = 31947.3
$ws.Range("M113").Value = -29777.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 3926.5
$ws.Range("I6").Value = 1000
$ws.Range("J6").Value = 4902
$ws.Range("K6").Value = 1000
$ws.Range("L6").Value = 4902
$ws.Range("M6").Value = -885
$ws.Range("N6").Value = -5132
